$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")

# Update test data values (normal load, cable capacitance etc.)
$ws.Range("D8").Value = 257
$ws.Range("E8").Value = 274
$ws.Range("D9").Value = 287
$ws.Range("E9").Value = 327

# Update the selected/active cell on the sheet view
$ws.Activate()
$ws.Range("F8").Select()
